# Edit script: 
#  1) Grow the "List Box" rectangle (Rectangle 1) taller.
#  2) Move the "Close" rectangle (Rectangle 17) down a bit.
#  3) Add a brand-new "Update" rectangle (Rectangle 10) next to the Close button.
#
# Because the shapes live inside floating drawing anchors (no plain body text),
# the most reliable way to make precise, faithful edits is to work on the
# document's WordprocessingML directly and push it back in with InsertXML.

$d = $word.ActiveDocument
$xml = $d.Content.WordOpenXML

# ---------------------------------------------------------------------------
# 1) Rectangle 1 ("List Box"): grow height  cy 3645752 -> 4212911
#    and effectExtent bottom   b  12065    -> 16510
#    (both the DrawingML wp:extent/a:ext sizes and the VML fallback height)
# ---------------------------------------------------------------------------
$xml = $xml.Replace(
    '<wp:extent cx="1457960" cy="3645752"/><wp:effectExtent l="0" t="0" r="15240" b="12065"/>',
    '<wp:extent cx="1457960" cy="4212911"/><wp:effectExtent l="0" t="0" r="15240" b="16510"/>')

$xml = $xml.Replace(
    '<a:off x="0" y="0"/><a:ext cx="1457960" cy="3645752"/>',
    '<a:off x="0" y="0"/><a:ext cx="1457960" cy="4212911"/>')

$xml = $xml.Replace(
    'width:114.8pt;height:287.05pt;',
    'width:114.8pt;height:331.75pt;')

# ---------------------------------------------------------------------------
# 2) Rectangle 17 ("Close"): reposition
#    positionH posOffset 3242310 -> 3241843
#    positionV posOffset 3300280 -> 3867110
#    VML fallback: margin-left:255.3pt -> 255.25pt, margin-top:259.85pt -> 304.5pt
# ---------------------------------------------------------------------------
$xml = $xml.Replace(
    '<wp:posOffset>3242310</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>3300280</wp:posOffset>',
    '<wp:posOffset>3241843</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>3867110</wp:posOffset>')

$xml = $xml.Replace(
    'margin-left:255.3pt;margin-top:259.85pt;width:105.7pt;height:33.7pt;',
    'margin-left:255.25pt;margin-top:304.5pt;width:105.7pt;height:33.7pt;')

# ---------------------------------------------------------------------------
# 3) Insert a brand new "Update" rectangle (Rectangle 10) right after the
#    "Close" rectangle (Rectangle 17) run.
# ---------------------------------------------------------------------------
$newShapeRun = '<w:r><w:rPr><w:noProof/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251687936" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="3A47B0F2" wp14:editId="6B5C01FD"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>3240911</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>3298785</wp:posOffset></wp:positionV><wp:extent cx="1330526" cy="427990"/><wp:effectExtent l="12700" t="12700" r="15875" b="16510"/><wp:wrapNone/><wp:docPr id="10" name="Rectangle 10"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="1330526" cy="427990"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></wps:spPr><wps:style><a:lnRef idx="3"><a:schemeClr val="lt1"/></a:lnRef><a:fillRef idx="1"><a:schemeClr val="accent3"/></a:fillRef><a:effectRef idx="1"><a:schemeClr val="accent3"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="lt1"/></a:fontRef></wps:style><wps:txbx><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Update</w:t></w:r></w:p></w:txbxContent></wps:txbx><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="ctr" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="margin"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:rect w14:anchorId="3A47B0F2" id="Rectangle 10" o:spid="_x0000_s1028" style="position:absolute;margin-left:255.2pt;margin-top:259.75pt;width:104.75pt;height:33.7pt;z-index:251687936;visibility:visible;mso-wrap-style:square;mso-width-percent:0;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;mso-width-percent:0;mso-width-relative:margin;v-text-anchor:middle" fillcolor="#a5a5a5 [3206]" strokecolor="white [3201]" strokeweight="1.5pt"><v:textbox><w:txbxContent><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Update</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect></w:pict></mc:Fallback></mc:AlternateContent></w:r>'

$closeRunEnd = '<v:rect w14:anchorId="31FC058E" id="Rectangle 17"'
$idx = $xml.IndexOf($closeRunEnd)
if ($idx -lt 0) {
    throw "Could not locate Rectangle 17 VML fallback to anchor the insertion"
}
$runClose = '</v:rect></w:pict></mc:Fallback></mc:AlternateContent></w:r>'
$closeIdx = $xml.IndexOf($runClose, $idx)
if ($closeIdx -lt 0) {
    throw "Could not locate end of Rectangle 17 run"
}
$insertAt = $closeIdx + $runClose.Length
$xml = $xml.Substring(0, $insertAt) + $newShapeRun + $xml.Substring($insertAt)

# ---------------------------------------------------------------------------
# Push the modified package XML back into the document.
# ---------------------------------------------------------------------------
$d.Content.InsertXML($xml)
Write-Host "Edit applied. Final XML length:" $xml.Length
